# Updated via Streamlit Approval System
# Appends three new pending-approval rows (84-86) to the sheet, mirroring
# the existing row layout (columns A:AO).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 84 : WGP011 ----
$ws.Cells.Item(84, 1).Value  = "WGP011"
$ws.Cells.Item(84, 2).Value  = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(84, 3).Value  = "13-02-2026"
$ws.Cells.Item(84, 4).Value  = 286962
$ws.Cells.Item(84, 5).Value  = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(84, 6).Value  = 34413429360
$ws.Cells.Item(84, 7).Value  = "NEFT"
$ws.Cells.Item(84, 8).Value  = "SBIN0003229"
$ws.Cells.Item(84, 9).Value  = "AAAFW8862C"
$ws.Cells.Item(84, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(84, 11).Value = "SHREE BALAJI ELECTRICAL"
$ws.Cells.Item(84, 12).Value = "d892dbf3-8741-44f1-ba7f-5a584bc7f350"
$ws.Cells.Item(84, 13).Value = "ACC-125006695576"
$ws.Cells.Item(84, 14).Value = "CNRB0017203"
$ws.Cells.Item(84, 21).Value = "pending"
$ws.Cells.Item(84, 22).Value = 105987
$ws.Cells.Item(84, 24).Value = "Being electric consumables purchased RPA_ID : abb7bb472c"
$ws.Cells.Item(84, 25).Value = "ONGC Electrical"
$ws.Cells.Item(84, 26).Value = "SITE EXPENSE"
$ws.Cells.Item(84, 27).Value = "midhuncraju12@gmail.com"
$ws.Cells.Item(84, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(84, 29).Value = 0
$ws.Cells.Item(84, 30).Value = 0
$ws.Cells.Item(84, 31).Value = 0

# ---- Row 85 : WGG 02 ----
$ws.Cells.Item(85, 1).Value  = "WGG 02"
$ws.Cells.Item(85, 2).Value  = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(85, 3).Value  = "13-02-2026"
$ws.Cells.Item(85, 4).Value  = 286962
$ws.Cells.Item(85, 5).Value  = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(85, 6).Value  = 34413429360
$ws.Cells.Item(85, 7).Value  = "NEFT"
$ws.Cells.Item(85, 8).Value  = "SBIN0003229"
$ws.Cells.Item(85, 9).Value  = "AAAFW8862C"
$ws.Cells.Item(85, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(85, 12).Value = "788a71f1-06f3-4161-8e96-7dbdaaa092ca"
$ws.Cells.Item(85, 21).Value = "pending"
$ws.Cells.Item(85, 22).Value = 1470
$ws.Cells.Item(85, 24).Value = "Being IOCL Willington switch and core cutting charges RPA_ID : 183a5be1f1"
$ws.Cells.Item(85, 25).Value = "IOCL Willington"
$ws.Cells.Item(85, 26).Value = "SITE EXPENSE"
$ws.Cells.Item(85, 27).Value = "midhuncraju12@gmail.com"
$ws.Cells.Item(85, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(85, 29).Value = 0
$ws.Cells.Item(85, 30).Value = 0
$ws.Cells.Item(85, 31).Value = 0

# ---- Row 86 : WGP008 ----
$ws.Cells.Item(86, 1).Value  = "WGP008"
$ws.Cells.Item(86, 2).Value  = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(86, 3).Value  = "13-02-2026"
$ws.Cells.Item(86, 4).Value  = 286962
$ws.Cells.Item(86, 5).Value  = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(86, 6).Value  = 34413429360
$ws.Cells.Item(86, 7).Value  = "DCR"
$ws.Cells.Item(86, 8).Value  = "SBIN0003229"
$ws.Cells.Item(86, 9).Value  = "AAAFW8862C"
$ws.Cells.Item(86, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(86, 11).Value = "GAYATHRI ELECTRICALS"
$ws.Cells.Item(86, 12).Value = "77e4cafa-fced-4e92-af02-1d695f6c561d"
$ws.Cells.Item(86, 13).Value = "ACC-39177475703"
$ws.Cells.Item(86, 14).Value = "SBIN0000512"
$ws.Cells.Item(86, 21).Value = "pending"
$ws.Cells.Item(86, 22).Value = 580870
$ws.Cells.Item(86, 24).Value = "Being material purchase RPA_ID : 32aed1f5ef"
$ws.Cells.Item(86, 25).Value = "ONGC Electrical"
$ws.Cells.Item(86, 26).Value = "SITE EXPENSES"
$ws.Cells.Item(86, 27).Value = "midhuncraju12@gmail.com"
$ws.Cells.Item(86, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(86, 29).Value = 0
$ws.Cells.Item(86, 30).Value = 0
$ws.Cells.Item(86, 31).Value = 0
